$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns with the latest
# scrape values. The D-column price strings are written via a quoted
# formula ("=""<value>""") and then collapsed back to a literal with
# Copy + PasteSpecial so Excel doesn't auto-convert decimal-looking
# text (e.g. "577.86") into a Number cell; the E-column percentages
# keep their original "  +/-x.xx%  " padded text layout and are never
# number-like (leading/trailing spaces), so a plain .Value assignment
# is safe for them.

$ws.Range("D2").Formula = "=""69.409.55"""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial()
$ws.Range("E2").Value = "  -2.20%  "

$ws.Range("D3").Formula = "=""3.476.85"""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial()
$ws.Range("E3").Value = "  -4.38%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Formula = "=""577.86"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial()
$ws.Range("E5").Value = "  -4.51%  "

$ws.Range("D6").Formula = "=""192.77"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial()
$ws.Range("E6").Value = "  -3.44%  "

$ws.Range("E7").Value = "  -3.13%  "

$ws.Range("D8").Formula = "=""3.466.27"""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial()
$ws.Range("E8").Value = "  -4.38%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("E10").Value = "  -8.12%  "

$ws.Range("D11").Formula = "=""0.618"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial()
$ws.Range("E11").Value = "  -4.65%  "

$ws.Range("D12").Formula = "=""51.39"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial()
$ws.Range("E12").Value = "  -4.78%  "

$ws.Range("E13").Value = "  -7.11%  "

$ws.Range("E14").Value = "  -4.64%  "

$ws.Range("D15").Formula = "=""4.040.51"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial()
$ws.Range("E15").Value = "  -4.12%  "

$ws.Range("D16").Formula = "=""642.34"""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial()
$ws.Range("E16").Value = "  -0.47%  "

$ws.Range("D17").Formula = "=""69.185.31"""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial()
$ws.Range("E17").Value = "  -2.62%  "

$ws.Range("D18").Formula = "=""3.467.44"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial()
$ws.Range("E18").Value = "  -4.37%  "

$ws.Range("D19").Formula = "=""12.34"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial()
$ws.Range("E19").Value = "  -5.10%  "

$ws.Range("D21").Formula = "=""18.19"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial()
$ws.Range("E21").Value = "  -4.70%  "

$ws.Range("D22").Formula = "=""0.944"""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial()
$ws.Range("E22").Value = "  -5.68%  "

$ws.Range("D23").Formula = "=""17.76"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial()
$ws.Range("E23").Value = "  -4.89%  "

$ws.Range("E24").Value = "  -1.15%  "

$ws.Range("D25").Formula = "=""98.99"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial()
$ws.Range("E25").Value = "  -5.22%  "

$ws.Range("E26").Value = "  -7.75%  "

$ws.Range("D27").Formula = "=""2.87"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial()
$ws.Range("E27").Value = "  -4.64%  "

$ws.Range("D28").Formula = "=""9.94"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial()
$ws.Range("E28").Value = "  -5.09%  "

$ws.Range("E29").Value = "  -4.54%  "

$ws.Range("D30").Formula = "=""32.44"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial()
$ws.Range("E30").Value = "  -4.70%  "

$ws.Range("D31").Formula = "=""4.30"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial()
$ws.Range("E31").Value = "  -9.91%  "

$ws.Range("D32").Formula = "=""6.73"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial()
$ws.Range("E32").Value = "  -6.56%  "

$ws.Range("D33").Formula = "=""11.62"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial()
$ws.Range("E33").Value = "  -5.17%  "

$ws.Range("E34").Value = "  -5.63%  "

$ws.Range("D35").Formula = "=""60.90"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial()
$ws.Range("E35").Value = "  -3.93%  "

$ws.Range("D36").Formula = "=""3.732.01"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial()
$ws.Range("E36").Value = "  -7.53%  "

$ws.Range("D37").Formula = "=""523.11"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial()
$ws.Range("E37").Value = "  +2.56%  "

$ws.Range("E38").Value = "  +0.04%  "

$ws.Range("D39").Formula = "=""0.0₃0793"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial()
$ws.Range("E39").Value = "  -9.88%  "

$ws.Range("D40").Formula = "=""2.95"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial()
$ws.Range("E40").Value = "  -3.33%  "

$ws.Range("D41").Formula = "=""3.49"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial()
$ws.Range("E41").Value = "  -1.66%  "

$ws.Range("D42").Formula = "=""0.372"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial()
$ws.Range("E42").Value = "  -4.73%  "

$ws.Range("E43").Value = "  -2.02%  "

$ws.Range("D44").Formula = "=""3.51"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial()
$ws.Range("E44").Value = "  +69.70%  "

$ws.Range("D45").Formula = "=""34.28"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial()
$ws.Range("E45").Value = "  -6.79%  "

$ws.Range("D46").Formula = "=""0.0443"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial()
$ws.Range("E46").Value = "  -4.14%  "

$ws.Range("D47").Formula = "=""3.35"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial()
$ws.Range("E47").Value = "  -4.29%  "

$ws.Range("E48").Value = "  -7.56%  "

$ws.Range("E49").Value = "  -4.50%  "

$ws.Range("E50").Value = "  -0.39%  "

$ws.Range("D51").Formula = "=""8.16"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial()
$ws.Range("E51").Value = "  -6.08%  "

$excel.CutCopyMode = $false